$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Green 5")
$ws.Name = "Green 5 (broken)"
$ws.Range("AJ4").Value = "broken"
$ws.Range("AJ4").Select()
